$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert the new "Make Sagan AI with nix" paragraph right after the
#    AI section's last bullet ("Arius 2nd miner doesn't mine sometimes
#    as well"), before the blank paragraph / "Battle:" heading.
# ------------------------------------------------------------------
$aiLast = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*mine sometimes as well*") {
        $aiLast = $d.Paragraphs.Item($i)
        break
    }
}

$aiLast.Range.InsertParagraphAfter()
$newPara = $aiLast.Next()
$newPara.Range.Text = "Make Sagan AI with nix"

# Position right after "Make Sagan " (11 characters) inside the new
# paragraph - this is where the _GoBack bookmark needs to live, split
# across the two runs.
$splitPoint = $newPara.Range.Start + 11

# ------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from wherever it currently sits
#    (end of "Have clumped units spread out a little bit") to the new
#    split point.
# ------------------------------------------------------------------
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 3) Relocate the <w:lastRenderedPageBreak/> marker: it currently sits
#    on the "Need to add discovered new unit..." run, but after the
#    new content is added earlier in the doc it should render on the
#    previous paragraph ("Fix splash reward particles not working")
#    instead. Rebuild both paragraphs verbatim via InsertXML so only
#    the break marker's position changes.
# ------------------------------------------------------------------
$pA = $null
$pB = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Fix splash reward particles not working*") {
        $pA = $d.Paragraphs.Item($i)
        $pB = $pA.Next()
        break
    }
}

$combined = $d.Range($pA.Range.Start, $pB.Range.End)

$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00313C47" w:rsidRDefault="00313C47" w:rsidP="00313C47"><w:r><w:lastRenderedPageBreak/><w:t>Fix splash reward particles not working</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00313C47" w:rsidRDefault="00313C47" w:rsidP="00313C47"><w:r><w:t xml:space="preserve">Need to add discovered new unit to reward splash screen  </w:t></w:r></w:p>'

[void]$combined.InsertXML($xmlFrag)
